# Seattle-Ready SPL admin-strings translation workbook update
# Adds the new "What to Expect at This Location" row (English source string
# plus its Spanish translation) to the bottom of the translation sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new row of translated strings --------------------------------
$rowA = $ws.Range("A93")
$rowB = $ws.Range("B93")

$rowA.Value = "What to Expect at This Location"
$rowB.Value = "Qué esperar en esta ubicación"

# Reset inherited column formatting (columns A:B carry a wrap-text style)
# before applying the per-cell fonts, so the new cells get a style record
# of their own instead of dragging the wrap/alignment settings along.
$rowA.Style = "Normal"
$rowB.Style = "Normal"

$rowA.Font.Name = "Times New Roman"
$rowB.Font.Name = "Calibri"

# --- Update the view so the new row is visible/selected -------------------
$ws.Range("A93:B93").Select()
